$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 22.46000000000007
$ws.Range("H2").Value = 0.0001097546089009471
$ws.Range("I2").Value = 0.0001097546089009471
$ws.Range("L2").Value = 51.3595395211851
$ws.Range("M2").Value = "[26.351031116638424, 76.36804792573177]"
$ws.Range("N2").Value = 0.0001522200079715486
$ws.Range("O2").Value = 0.0001522200079715486
$ws.Range("P2").Value = 1.490605523324886
$ws.Range("Q2").Value = "[0.8993948938205003, 2.081816152829272]"
$ws.Range("R2").Value = [double]"7.089210066091312e-06"
$ws.Range("S2").Value = [double]"7.089210066091312e-06"
$ws.Range("T2").Value = 64.81571568110809
$ws.Range("U2").Value = "[49.81730060115501, 79.81413076106116]"
$ws.Range("V2").Value = [double]"3.33377769834442e-11"
$ws.Range("W2").Value = [double]"3.33377769834442e-11"
$ws.Range("X2").Value = 17.13165165165171
$ws.Range("Y2").Value = 15.01829829829835
$ws.Range("Z2").Value = 19.24500500500507

# Row 3
$ws.Range("F3").Value = 22.46000000000007
$ws.Range("H3").Value = 0.02405451339997511
$ws.Range("I3").Value = 0.02405451339997511
$ws.Range("L3").Value = 29.16022859184519
$ws.Range("M3").Value = "[2.291902773396039, 56.02855441029434]"
$ws.Range("N3").Value = 0.0340656162125923
$ws.Range("O3").Value = 0.0340656162125923
$ws.Range("P3").Value = 1.842816111114733
$ws.Range("Q3").Value = "[0.6352369529781159, 3.050395269251351]"
$ws.Range("R3").Value = 0.003586295186198019
$ws.Range("S3").Value = 0.003586295186198019
$ws.Range("T3").Value = 56.92540049893095
$ws.Range("U3").Value = "[42.62456117249222, 71.22623982536969]"
$ws.Range("V3").Value = [double]"3.253326497087983e-10"
$ws.Range("W3").Value = [double]"3.253326497087983e-10"
$ws.Range("X3").Value = 15.87263263263268
$ws.Range("Y3").Value = 11.55599599599603
$ws.Range("Z3").Value = 20.18926926926934

# Row 4
$ws.Range("F4").Value = 22.46000000000007
$ws.Range("H4").Value = [double]"3.13237495686014e-05"
$ws.Range("I4").Value = [double]"3.13237495686014e-05"
$ws.Range("L4").Value = 48.23120380098898
$ws.Range("M4").Value = "[21.7040612407, 74.75834636127796]"
$ws.Range("N4").Value = 0.0006558070537809879
$ws.Range("O4").Value = 0.0006558070537809879
$ws.Range("P4").Value = 2.421447791055196
$ws.Range("Q4").Value = "[1.9182898084982725, 2.924605773612119]"
$ws.Range("R4").Value = [double]"1.368904989362818e-12"
$ws.Range("S4").Value = [double]"1.368904989362818e-12"
$ws.Range("T4").Value = 59.20191066838621
$ws.Range("U4").Value = "[45.65280598364599, 72.75101535312643]"
$ws.Range("V4").Value = [double]"2.429345613563783e-11"
$ws.Range("W4").Value = [double]"2.429345613563783e-11"
$ws.Range("X4").Value = 13.80424424424429
$ws.Range("Y4").Value = 12.00564564564569
$ws.Range("Z4").Value = 15.60284284284289

# Row 5
$ws.Range("F5").Value = 22.46000000000007
$ws.Range("H5").Value = [double]"6.121321586394224e-06"
$ws.Range("I5").Value = [double]"6.121321586394224e-06"
$ws.Range("L5").Value = 43.3976678894764
$ws.Range("M5").Value = "[24.257125755082477, 62.53821002387033]"
$ws.Range("N5").Value = [double]"3.83438594315777e-05"
$ws.Range("O5").Value = [double]"3.83438594315777e-05"
$ws.Range("P5").Value = 2.974921571867812
$ws.Range("Q5").Value = "[2.446605690183042, 3.5032374535525816]"
$ws.Range("R5").Value = [double]"8.659739592076221e-15"
$ws.Range("S5").Value = [double]"8.659739592076221e-15"
$ws.Range("T5").Value = 57.16036169620677
$ws.Range("U5").Value = "[45.86159471067129, 68.45912868174224]"
$ws.Range("V5").Value = [double]"2.879918525877656e-13"
$ws.Range("W5").Value = [double]"2.879918525877656e-13"
$ws.Range("X5").Value = 11.82578578578582
$ws.Range("Y5").Value = 9.93725725725729
$ws.Range("Z5").Value = 13.71431431431436

# Row 6
$ws.Range("F6").Value = 22.46000000000007
$ws.Range("H6").Value = 0.01945163481364709
$ws.Range("I6").Value = 0.01945163481364709
$ws.Range("L6").Value = 27.19800559277141
$ws.Range("M6").Value = "[2.844485381484084, 51.551525804058734]"
$ws.Range("N6").Value = 0.02942959339717999
$ws.Range("O6").Value = 0.02942959339717999
$ws.Range("P6").Value = -2.905737349266235
$ws.Range("Q6").Value = "[-4.138474406530699, -1.6730002920017717]"
$ws.Range("R6").Value = [double]"2.121576289471072e-05"
$ws.Range("S6").Value = [double]"2.121576289471072e-05"
$ws.Range("T6").Value = 54.63288848940283
$ws.Range("U6").Value = "[41.134330555364826, 68.13144642344083]"
$ws.Range("V6").Value = [double]"2.075435379111923e-10"
$ws.Range("W6").Value = [double]"2.075435379111923e-10"
$ws.Range("X6").Value = 10.38690690690694
$ws.Range("Y6").Value = 5.980340340340359
$ws.Range("Z6").Value = 14.79347347347352

# Row 7
$ws.Range("F7").Value = 22.46000000000007
$ws.Range("H7").Value = [double]"1.141585476260865e-06"
$ws.Range("I7").Value = [double]"1.141585476260865e-06"
$ws.Range("L7").Value = 54.2763601707308
$ws.Range("M7").Value = "[31.012033697429715, 77.54068664403188]"
$ws.Range("N7").Value = [double]"2.488435445791026e-05"
$ws.Range("O7").Value = [double]"2.488435445791026e-05"
$ws.Range("P7").Value = -2.805105752754851
$ws.Range("Q7").Value = "[-3.2956847857478513, -2.31452671976185]"
$ws.Range("R7").Value = [double]"5.329070518200751e-15"
$ws.Range("S7").Value = [double]"5.329070518200751e-15"
$ws.Range("T7").Value = 65.62917440210805
$ws.Range("U7").Value = "[52.818201093497365, 78.44014771071873]"
$ws.Range("V7").Value = [double]"1.931788062847772e-13"
$ws.Range("W7").Value = [double]"1.931788062847772e-13"
$ws.Range("X7").Value = 10.02718718718722
$ws.Range("Y7").Value = 8.273553553553583
$ws.Range("Z7").Value = 11.78082082082086

# Row 8
$ws.Range("F8").Value = 25.63000000000057
$ws.Range("H8").Value = 0.0002090287051502715
$ws.Range("I8").Value = 0.0002090287051502715
$ws.Range("L8").Value = 37.74762089259059
$ws.Range("M8").Value = "[19.246136902628812, 56.249104882552366]"
$ws.Range("N8").Value = 0.0001657479154211217
$ws.Range("O8").Value = 0.0001657479154211217
$ws.Range("P8").Value = 3.125868966634889
$ws.Range("Q8").Value = "[2.5220793875665803, 3.729658545703198]"
$ws.Range("R8").Value = [double]"1.381117442633695e-13"
$ws.Range("S8").Value = [double]"1.381117442633695e-13"
$ws.Range("T8").Value = 48.56907141658719
$ws.Range("U8").Value = "[37.06724100096404, 60.070901832210346]"
$ws.Range("V8").Value = [double]"6.419997866657923e-11"
$ws.Range("W8").Value = [double]"6.419997866657923e-11"
$ws.Range("X8").Value = 12.87913913913943
$ws.Range("Y8").Value = 10.41619619619643
$ws.Range("Z8").Value = 15.34208208208242

# Row 9
$ws.Range("F9").Value = 25.63000000000057
$ws.Range("H9").Value = 0.03928183839056287
$ws.Range("I9").Value = 0.03928183839056287
$ws.Range("L9").Value = 29.79526024661671
$ws.Range("M9").Value = "[2.6348522404648094, 56.95566825276862]"
$ws.Range("N9").Value = 0.03227170641516519
$ws.Range("O9").Value = 0.03227170641516519
$ws.Range("P9").Value = -3.056684744033312
$ws.Range("Q9").Value = "[-4.314579700425622, -1.7987897876410028]"
$ws.Range("R9").Value = [double]"1.306991498206678e-05"
$ws.Range("S9").Value = [double]"1.306991498206678e-05"
$ws.Range("T9").Value = 50.99773961505379
$ws.Range("U9").Value = "[35.44454079446382, 66.55093843564376]"
$ws.Range("V9").Value = [double]"3.946921989239627e-08"
$ws.Range("W9").Value = [double]"3.946921989239627e-08"
$ws.Range("X9").Value = 12.46864864864893
$ws.Range("Y9").Value = 7.337517517517683
$ws.Range("Z9").Value = 17.59977977978017

# Row 10
$ws.Range("B10").Value = 1
$ws.Range("F10").Value = 25.63000000000057
$ws.Range("H10").Value = 0.0006301998591198288
$ws.Range("I10").Value = 0.0006301998591198288
$ws.Range("L10").Value = 37.46242616282227
$ws.Range("M10").Value = "[15.791368139313263, 59.13348418633127]"
$ws.Range("N10").Value = 0.001119960324487135
$ws.Range("O10").Value = 0.001119960324487135
$ws.Range("P10").Value = 2.937184723176043
$ws.Range("Q10").Value = "[2.232763547596349, 3.6416058987557367]"
$ws.Range("R10").Value = [double]"9.145151302902832e-11"
$ws.Range("S10").Value = [double]"9.145151302902832e-11"
$ws.Range("T10").Value = 55.18359740633262
$ws.Range("U10").Value = "[42.65072867902526, 67.71646613363997]"
$ws.Range("V10").Value = [double]"1.946132144325929e-11"
$ws.Range("W10").Value = [double]"1.946132144325929e-11"
$ws.Range("X10").Value = 13.64880880880911
$ws.Range("Y10").Value = 10.77537537537561
$ws.Range("Z10").Value = 16.52224224224261

# Row 11
$ws.Range("F11").Value = 25.63000000000057
$ws.Range("H11").Value = 0.0002977775069856747
$ws.Range("I11").Value = 0.0002977775069856747
$ws.Range("L11").Value = 44.76471722542993
$ws.Range("M11").Value = "[16.605168997150102, 72.92426545370976]"
$ws.Range("N11").Value = 0.002507475860896502
$ws.Range("O11").Value = 0.002507475860896502
$ws.Range("P11").Value = 2.081816152829273
$ws.Range("Q11").Value = "[1.478026573760964, 2.6856057318975814]"
$ws.Range("R11").Value = [double]"1.231825685898968e-08"
$ws.Range("S11").Value = [double]"1.231825685898968e-08"
$ws.Range("T11").Value = 60.70110345930016
$ws.Range("U11").Value = "[45.802078569646426, 75.6001283489539]"
$ws.Range("V11").Value = [double]"1.732884946648028e-10"
$ws.Range("W11").Value = [double]"1.732884946648028e-10"
$ws.Range("X11").Value = 17.13797797797836
$ws.Range("Y11").Value = 14.67503503503536
$ws.Range("Z11").Value = 19.60092092092135

# Row 12
$ws.Range("F12").Value = 25.63000000000057
$ws.Range("H12").Value = 0.00208668504259657
$ws.Range("I12").Value = 0.00208668504259657
$ws.Range("L12").Value = 43.06811796255133
$ws.Range("M12").Value = "[13.24269256915693, 72.89354335594572]"
$ws.Range("N12").Value = 0.005624985652465941
$ws.Range("O12").Value = 0.005624985652465941
$ws.Range("P12").Value = 1.792500312859041
$ws.Range("Q12").Value = "[0.9497106920761933, 2.6352899336418893]"
$ws.Range("R12").Value = [double]"9.540116799744425e-05"
$ws.Range("S12").Value = [double]"9.540116799744425e-05"
$ws.Range("T12").Value = 56.47503150676456
$ws.Range("U12").Value = "[39.59249702518483, 73.35756598834429]"
$ws.Range("V12").Value = [double]"2.499546947198894e-08"
$ws.Range("W12").Value = [double]"2.499546947198894e-08"
$ws.Range("X12").Value = 18.31813813813854
$ws.Range("Y12").Value = 14.88028028028061
$ws.Range("Z12").Value = 21.75599599599648

# Row 13
$ws.Range("F13").Value = 25.63000000000057
$ws.Range("H13").Value = [double]"7.881273630716557e-05"
$ws.Range("I13").Value = [double]"7.881273630716557e-05"
$ws.Range("L13").Value = 45.80726281351437
$ws.Range("M13").Value = "[20.789400268588025, 70.82512535844072]"
$ws.Range("N13").Value = 0.0006069018876757593
$ws.Range("O13").Value = 0.0006069018876757593
$ws.Range("P13").Value = 1.213868632918579
$ws.Range("Q13").Value = "[0.6226580034141929, 1.8050792624229643]"
$ws.Range("R13").Value = 0.0001526951013641753
$ws.Range("S13").Value = 0.0001526951013641753
$ws.Range("T13").Value = 63.02265973336081
$ws.Range("U13").Value = "[49.38764227388151, 76.65767719284011]"
$ws.Range("V13").Value = [double]"4.658273766722232e-12"
$ws.Range("W13").Value = [double]"4.658273766722232e-12"
$ws.Range("X13").Value = 20.67845845845892
$ws.Range("Y13").Value = 18.26682682682723
$ws.Range("Z13").Value = 23.0900900900906
